$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.036300787868605
$ws.Range("D2").Value = 1.039486257509548
$ws.Range("E2").Value = 1.044854920479544
$ws.Range("F2").Value = 1.054943044905987
$ws.Range("I2").Value = 1.03750406556183
$ws.Range("J2").Value = 1.041409853827325
$ws.Range("K2").Value = 1.042271308535869
$ws.Range("L2").Value = 1.047624798620724
$ws.Range("M2").Value = 1.057684864806386
$ws.Range("N2").Value = 1.042888775271139

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.037185561383763
$ws.Range("D3").Value = 1.040150680955219
$ws.Range("E3").Value = 1.045714901049047
$ws.Range("F3").Value = 1.056052238088701
$ws.Range("I3").Value = 1.037700139899264
$ws.Range("J3").Value = 1.041938860580838
$ws.Range("K3").Value = 1.042746344899974
$ws.Range("L3").Value = 1.04829597473905
$ws.Range("M3").Value = 1.058606645936786
$ws.Range("N3").Value = 1.043418533274921

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.037758498825932
$ws.Range("D4").Value = 1.040580904628125
$ws.Range("E4").Value = 1.046272184856157
$ws.Range("F4").Value = 1.056771304409433
$ws.Range("I4").Value = 1.037826013803047
$ws.Range("J4").Value = 1.042280948434299
$ws.Range("K4").Value = 1.043053339620568
$ws.Range("L4").Value = 1.048730446826958
$ws.Range("M4").Value = 1.059203851114266
$ws.Range("N4").Value = 1.043761106932371

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.037999463463835
$ws.Range("D5").Value = 1.040761840388028
$ws.Range("E5").Value = 1.046506661716478
$ws.Range("F5").Value = 1.057073920590527
$ws.Range("I5").Value = 1.037878691406156
$ws.Range("J5").Value = 1.042424709754515
$ws.Range("K5").Value = 1.043182306886923
$ws.Range("L5").Value = 1.048913140060921
$ws.Range("M5").Value = 1.059455094938615
$ws.Range("N5").Value = 1.043905072410152

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038039928419086
$ws.Range("D6").Value = 1.040792224345252
$ws.Range("E6").Value = 1.046546042795668
$ws.Range("F6").Value = 1.057124749953201
$ws.Range("I6").Value = 1.037887522130373
$ws.Range("J6").Value = 1.042448844811097
$ws.Range("K6").Value = 1.04320395555717
$ws.Range("L6").Value = 1.04894381744307
$ws.Range("M6").Value = 1.059497290334607
$ws.Range("N6").Value = 1.043929241741282

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.037761718208905
$ws.Range("D7").Value = 1.040583322028431
$ws.Range("E7").Value = 1.04627531718383
$ws.Range("F7").Value = 1.056775346722884
$ws.Range("I7").Value = 1.037826718626638
$ws.Range("J7").Value = 1.042282869587255
$ws.Range("K7").Value = 1.043055063255415
$ws.Range("L7").Value = 1.048732887822412
$ws.Range("M7").Value = 1.059207207545426
$ws.Range("N7").Value = 1.043763030813584

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.036599711525193
$ws.Range("D8").Value = 1.039710740105378
$ws.Range("E8").Value = 1.045145384608738
$ws.Range("F8").Value = 1.055317622957431
$ws.Range("I8").Value = 1.037570536382943
$ws.Range("J8").Value = 1.041588678048191
$ws.Range("K8").Value = 1.042431928528408
$ws.Range("L8").Value = 1.047851588617595
$ws.Range("M8").Value = 1.057996228886971
$ws.Range("N8").Value = 1.043067853442911

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.034555450503612
$ws.Range("D9").Value = 1.038175479643083
$ws.Range("E9").Value = 1.043160627755306
$ws.Range("F9").Value = 1.052759270589346
$ws.Range("I9").Value = 1.037111481444725
$ws.Range("J9").Value = 1.040363820854391
$ws.Range("K9").Value = 1.041330977019868
$ws.Range("L9").Value = 1.046300022227499
$ws.Range("M9").Value = 1.055868129516777
$ws.Range("N9").Value = 1.041841256811402

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.033194925940619
$ws.Range("D10").Value = 1.037153631619897
$ws.Range("E10").Value = 1.041841787940188
$ws.Range("F10").Value = 1.051060720705059
$ws.Range("I10").Value = 1.036800348092571
$ws.Range("J10").Value = 1.039546228826786
$ws.Range("K10").Value = 1.040595113461808
$ws.Range("L10").Value = 1.045266643520442
$ws.Range("M10").Value = 1.054453358243677
$ws.Range("N10").Value = 1.041022503709338

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.032606367243291
$ws.Range("D11").Value = 1.036711571534881
$ws.Range("E11").Value = 1.04127175898505
$ws.Range("F11").Value = 1.050326907676914
$ws.Range("I11").Value = 1.036664421835066
$ws.Range("J11").Value = 1.039191972287152
$ws.Range("K11").Value = 1.040276039179043
$ws.Range("L11").Value = 1.044819428646358
$ws.Range("M11").Value = 1.053841699118587
$ws.Range("N11").Value = 1.040667744084785

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.03238783518029
$ws.Range("D12").Value = 1.036547433292909
$ws.Range("E12").Value = 1.041060181964107
$ws.Range("F12").Value = 1.050054588805302
$ws.Range("I12").Value = 1.036613752486783
$ws.Range("J12").Value = 1.039060351442369
$ws.Range("K12").Value = 1.040157455581376
$ws.Range("L12").Value = 1.044653350913284
$ws.Range("M12").Value = 1.053614644830838
$ws.Range("N12").Value = 1.040535936323305

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.032434707192651
$ws.Range("D13").Value = 1.036582638649216
$ws.Range("E13").Value = 1.041105558823376
$ws.Range("F13").Value = 1.050112990709512
$ws.Range("I13").Value = 1.036624629383815
$ws.Range("J13").Value = 1.039088586092672
$ws.Range("K13").Value = 1.040182895100431
$ws.Range("L13").Value = 1.044688973425151
$ws.Range("M13").Value = 1.05366334225267
$ws.Range("N13").Value = 1.04056421107005

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.032588301577661
$ws.Range("D14").Value = 1.036698002532191
$ws.Range("E14").Value = 1.041254266750145
$ws.Range("F14").Value = 1.050304392568206
$ws.Range("I14").Value = 1.036660237165641
$ws.Range("J14").Value = 1.039181093159182
$ws.Range("K14").Value = 1.040266238346309
$ws.Range("L14").Value = 1.044805699832649
$ws.Range("M14").Value = 1.05382292781444
$ws.Range("N14").Value = 1.040656849507206

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.032682947422915
$ws.Range("D15").Value = 1.036769090360099
$ws.Range("E15").Value = 1.041345911468541
$ws.Range("F15").Value = 1.050422354988168
$ws.Range("I15").Value = 1.036682152421505
$ws.Range("J15").Value = 1.039238085313043
$ws.Range("K15").Value = 1.040317580260134
$ws.Range("L15").Value = 1.044877623860984
$ws.Range("M15").Value = 1.053921272744017
$ws.Range("N15").Value = 1.040713922596462

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.033233998480666
$ws.Range("D16").Value = 1.037182978385963
$ws.Range("E16").Value = 1.041879640848468
$ws.Range("F16").Value = 1.051109456736725
$ws.Range("I16").Value = 1.036809343759592
$ws.Range("J16").Value = 1.039569734820345
$ws.Range("K16").Value = 1.040616280173674
$ws.Range("L16").Value = 1.045296328943172
$ws.Range("M16").Value = 1.05449397206482
$ws.Range("N16").Value = 1.041046043084105

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.03357980826003
$ws.Range("D17").Value = 1.037442709346518
$ws.Range("E17").Value = 1.042214713812893
$ws.Range("F17").Value = 1.051540905409684
$ws.Range("I17").Value = 1.036888805692644
$ws.Range("J17").Value = 1.039777707858085
$ws.Range("K17").Value = 1.040803529492541
$ws.Range("L17").Value = 1.045559037847947
$ws.Range("M17").Value = 1.054853465235303
$ws.Range("N17").Value = 1.041254311467413

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.033781566962483
$ws.Range("D18").Value = 1.037594245180915
$ws.Range("E18").Value = 1.042410256212463
$ws.Range("F18").Value = 1.051792723210216
$ws.Range("I18").Value = 1.036935038341921
$ws.Range("J18").Value = 1.039898992391773
$ws.Range("K18").Value = 1.04091270628407
$ws.Range("L18").Value = 1.045712295072871
$ws.Range("M18").Value = 1.055063242744461
$ws.Range("N18").Value = 1.041375768239051

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.033850370542735
$ws.Range("D19").Value = 1.037645921582958
$ws.Range("E19").Value = 1.042476948042378
$ws.Range("F19").Value = 1.051878613865239
$ws.Range("I19").Value = 1.036950782763326
$ws.Range("J19").Value = 1.039940343418296
$ws.Range("K19").Value = 1.040949925498539
$ws.Range("L19").Value = 1.0457645557932
$ws.Range("M19").Value = 1.05513478689713
$ws.Range("N19").Value = 1.041417177988775

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.033542700580116
$ws.Range("D20").Value = 1.037414838611734
$ws.Range("E20").Value = 1.042178753301471
$ws.Range("F20").Value = 1.051494598380665
$ws.Range("I20").Value = 1.036880292189605
$ws.Range("J20").Value = 1.039755396652379
$ws.Range("K20").Value = 1.040783443808682
$ws.Range("L20").Value = 1.045530849224093
$ws.Range("M20").Value = 1.054814885577289
$ws.Range("N20").Value = 1.041231968577237

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.032543069508607
$ws.Range("D21").Value = 1.036664029005675
$ws.Range("E21").Value = 1.041210471624945
$ws.Range("F21").Value = 1.050248022527838
$ws.Range("I21").Value = 1.036649756533911
$ws.Range("J21").Value = 1.039153853070162
$ws.Range("K21").Value = 1.04024169762971
$ws.Range("L21").Value = 1.044771325764031
$ws.Range("M21").Value = 1.053775929884543
$ws.Range("N21").Value = 1.040629570734134

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.031915053199355
$ws.Range("D22").Value = 1.036192327846222
$ws.Range("E22").Value = 1.040602584138794
$ws.Range("F22").Value = 1.04946570908541
$ws.Range("I22").Value = 1.036503766809198
$ws.Range("J22").Value = 1.03877544137944
$ws.Range("K22").Value = 1.039900703550288
$ws.Range("L22").Value = 1.044294002312695
$ws.Range("M22").Value = 1.053123525895736
$ws.Range("N22").Value = 1.040250621655407

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.032247929578746
$ws.Range("D23").Value = 1.036442350704021
$ws.Range("E23").Value = 1.04092475015332
$ws.Range("F23").Value = 1.049880289673601
$ws.Range("I23").Value = 1.036581257395556
$ws.Range("J23").Value = 1.038976062972563
$ws.Range("K23").Value = 1.040081506348495
$ws.Range("L23").Value = 1.044547019407942
$ws.Range("M23").Value = 1.053469298710487
$ws.Range("N23").Value = 1.040451528154205

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.033559467777529
$ws.Range("D24").Value = 1.037427432076268
$ws.Range("E24").Value = 1.042195001999795
$ws.Range("F24").Value = 1.051515522036143
$ws.Range("I24").Value = 1.036884139434108
$ws.Range("J24").Value = 1.039765478195846
$ws.Range("K24").Value = 1.040792519796268
$ws.Range("L24").Value = 1.045543586377122
$ws.Range("M24").Value = 1.054832317782304
$ws.Range("N24").Value = 1.041242064437652

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.035083537764911
$ws.Range("D25").Value = 1.038572095121158
$ws.Range("E25").Value = 1.043672977132548
$ws.Range("F25").Value = 1.053419433690068
$ws.Range("I25").Value = 1.037231058880745
$ws.Range("J25").Value = 1.040680659874693
$ws.Range("K25").Value = 1.041615937906686
$ws.Range("L25").Value = 1.046700967723134
$ws.Range("M25").Value = 1.056417600522256
$ws.Range("N25").Value = 1.042158545779455
